$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (stock alert records) appended starting at row 653, columns A:O,
# matching the existing table layout (header in row 1, data from row 2).
$newRows = @(
    @("BATAINDIA.NS",39417,65,38869,47,59.47750667254579,39295,61,59.69378497807872,"Low",0.01544845039520989,58.75142950397092,3,1,2),
    @("MANAPPURAM.NS",42917,84,42644,75,90.23010662436555,42795,80,90.57228138015006,"High",0.06843495115690246,85.09748528759786,3,1,2),
    @("GODREJPROP.NS",43344,103,43101,95,911.8499755859375,43221,99,920,"High",2.037506103515625,718.2868957519531,3,1,2),
    @("GMRAIRPORT.NS",42370,112,42217,107,8.829098701477051,42248,108,8.829098701477051,"Low",0,8.829098701477051,3,1,2),
    @("GMRAIRPORT.NS",42522,117,42217,107,8.829098701477051,42401,113,8.784051895141602,"Low",-0.007507801055908203,9.632433414459229,3,1,2),
    @("GMRAIRPORT.NS",42522,117,42248,108,8.829098701477051,42401,113,8.784051895141602,"Low",-0.009009361267089844,9.802109718322754,3,1,2),
    @("BATAINDIA.NS",42705,173,42217,157,581.4596786134197,42583,169,577.2650448407002,"High",-0.3495528143932916,636.3394704731664,3,1,2),
    @("MANAPPURAM.NS",44348,131,44013,120,170.1301489508114,44228,127,168.705578165982,"High",-0.2035101121184937,194.5513624050307,3,1,2),
    @("LICHSGFIN.NS",42461,165,42217,157,431.6130438757089,42339,161,428.005344062932,"High",-0.9019249531942393,573.2152615272045,3,1,2),
    @("LICHSGFIN.NS",42522,167,42156,155,320.6515052936053,42401,163,323.2496289079498,"Low",0.3247654517930698,270.3128602656794,3,1,2),
    @("GMRAIRPORT.NS",43983,165,43344,144,19.36996269226074,43525,150,19.23482322692871,"High",-0.02252324422200521,22.61330986022949,3,1,2),
    @("GMRAIRPORT.NS",44013,166,43617,153,12.65804481506348,43891,162,12.70309066772461,"Low",0.005005094740125868,11.89226531982422,3,1,2)

)

$startRow = 653
$dateCols = @(2, 4, 7)  # B, D, G columns (1-based) hold Excel date serials

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value2 = $rowData[$c - 1]
        if ($dateCols -contains $c) {
            $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
        }
    }
}
